# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new market day 2022-05-06) right after
# the existing row 6, pushing the rest of the data set down by two rows
# (old rows 7-35 become rows 9-37) and extending the used range to T37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 7 (shifts rows 7.. down to 9..)
$ws.Rows.Item(7).EntireRow.Insert()
$ws.Rows.Item(7).EntireRow.Insert()

# New row 7: Higo, Primera, 100 bandejas, 15000/15000/15000, 2022-05-06
$ws.Cells.Item(7,1).Value  = 6
$ws.Cells.Item(7,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(7,3).Value  = "Metropolitana"
$ws.Cells.Item(7,4).Value  = 44687
$ws.Cells.Item(7,5).Value  = 13
$ws.Cells.Item(7,6).Value  = "Fruta"
$ws.Cells.Item(7,7).Value  = 100101
$ws.Cells.Item(7,8).Value  = "Berries"
$ws.Cells.Item(7,9).Value  = 100101006
$ws.Cells.Item(7,10).Value = "Higo"
$ws.Cells.Item(7,11).Value = "Sin especificar"
$ws.Cells.Item(7,12).Value = "Primera"
$ws.Cells.Item(7,13).Value = 100
$ws.Cells.Item(7,14).Value = 15000
$ws.Cells.Item(7,15).Value = 15000
$ws.Cells.Item(7,16).Value = 15000
$ws.Cells.Item(7,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(7,18).Value = "Región Metropolitana"
$ws.Cells.Item(7,19).Value = 2143
$ws.Cells.Item(7,20).Value = 7

# New row 8: Higo, Segunda, 75 bandejas, 12000/12000/12000, 2022-05-06
$ws.Cells.Item(8,1).Value  = 6
$ws.Cells.Item(8,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(8,3).Value  = "Metropolitana"
$ws.Cells.Item(8,4).Value  = 44687
$ws.Cells.Item(8,5).Value  = 13
$ws.Cells.Item(8,6).Value  = "Fruta"
$ws.Cells.Item(8,7).Value  = 100101
$ws.Cells.Item(8,8).Value  = "Berries"
$ws.Cells.Item(8,9).Value  = 100101006
$ws.Cells.Item(8,10).Value = "Higo"
$ws.Cells.Item(8,11).Value = "Sin especificar"
$ws.Cells.Item(8,12).Value = "Segunda"
$ws.Cells.Item(8,13).Value = 75
$ws.Cells.Item(8,14).Value = 12000
$ws.Cells.Item(8,15).Value = 12000
$ws.Cells.Item(8,16).Value = 12000
$ws.Cells.Item(8,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(8,18).Value = "Región Metropolitana"
$ws.Cells.Item(8,19).Value = 1714
$ws.Cells.Item(8,20).Value = 7
